$d = $word.ActiveDocument

# 1) Lecturer placeholder: Dfsdf -> Sdfsdf
$d.Content.Find.Execute("Dfsdf", $true, $false, $false, $false, $false, $true, 1, $false, "Sdfsdf", 2) | Out-Null

# 2) Room placeholder: SFSDF -> SDFSDF
$d.Content.Find.Execute("SFSDF", $true, $false, $false, $false, $false, $true, 1, $false, "SDFSDF", 2) | Out-Null

# 3) Email placeholder: sdfds@sfdfds.sdfsdf -> sdfsdf@sdfsdf.sdfs
$d.Content.Find.Execute("sdfds@sfdfds.sdfsdf", $true, $false, $false, $false, $false, $true, 1, $false, "sdfsdf@sdfsdf.sdfs", 2) | Out-Null

# 4) The paragraph that holds the "Develop workplace soft-skills" text becomes the
#    "++For BA courses, list all learning outcomes ...++" instruction note, written
#    as several highlighted/bold Arial runs. Do this BEFORE the text-outcome swap
#    below so the "Develop workplace soft-skills" text is still unique in the
#    document when we search for it.
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "*Develop workplace soft-skills*") {
        $target = $cand
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the 'Develop workplace soft-skills' paragraph"
}
$r = $target.Range
$contentRange = $d.Range($r.Start, $r.End)
$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">++For </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>BA</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> courses,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>list all learning outcomes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> from the course descriptor, as well as Knowledge Gains and Personal Skills Development</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>outcomes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>.+</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/><w:highlight w:val="yellow"/></w:rPr><w:t>+</w:t></w:r>
</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange.InsertXML($frag)

# 5) Swap the text of the "text file operation" learning-outcome paragraph to the
#    "Develop workplace soft-skills" text (the old wording of that bullet is removed).
$oldOutcome = "Demonstrate an understanding of text file operation (reading/writing) and develop the appropriate program code for such operation including exception handling and data validation."
$newOutcome = "Develop workplace soft-skills including working in groups, writing formal reports, carrying out individual research and/or delivering oral presentations"
$d.Content.Find.Execute($oldOutcome, $true, $false, $false, $false, $false, $true, 1, $false, $newOutcome, 2) | Out-Null
